$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.307.68'
$ws.Range('E2').Value = '  +5.12%  '
$ws.Range('D3').Value = '3.506.50'
$ws.Range('E3').Value = '  +2.66%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''418.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('D6').Value = '''132.70'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.43%  '
$ws.Range('D7').Value = '''0.654'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.72%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '''0.777'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.91%  '
$ws.Range('E10').Value = '  +16.28%  '
$ws.Range('D11').Value = '''43.48'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').Value = '''0.0000266'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +22.17%  '
$ws.Range('D13').Value = '''10.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.89%  '
$ws.Range('D14').Value = '4.070.37'
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').Value = '''0.140'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '''20.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').Value = '3.504.81'
$ws.Range('E17').Value = '  +2.54%  '
$ws.Range('D18').Value = '''12.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('E19').Value = '  +3.37%  '
$ws.Range('D20').Value = '65.242.92'
$ws.Range('E20').Value = '  +5.02%  '
$ws.Range('D21').Value = '''454.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.88%  '
$ws.Range('D22').Value = '''90.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').Value = '''13.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('D25').Value = '''3.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.67%  '
$ws.Range('D26').Value = '''9.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.91%  '
$ws.Range('D27').Value = '''34.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.71%  '
$ws.Range('D28').Value = '''12.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.02%  '
$ws.Range('D29').Value = '''2.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.25%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '''0.118'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.29%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = '''7.45'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.82%  '
$ws.Range('E32').Value = '  -2.01%  '
$ws.Range('D33').Value = '''39.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('D34').Value = '''0.997'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').Value = '''57.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.04%  '
$ws.Range('D36').Value = '''0.0508'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.28%  '
$ws.Range('D37').Value = '0.0₃0738'
$ws.Range('E37').Value = '  +38.11%  '
$ws.Range('D38').Value = '''0.148'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.07%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '''0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''3.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').Value = '''4.55'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.23%  '
$ws.Range('E42').Value = '  +2.23%  '
$ws.Range('D43').Value = '''146.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('E45').Value = '  -2.96%  '
$ws.Range('E46').Value = '  -2.83%  '
$ws.Range('D47').Value = '''2.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').Value = '''15.88'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('D49').Value = '''0.146'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.26%  '
$ws.Range('D50').Value = '''2.58'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.63%  '
$ws.Range('D51').Value = '''21.81'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.91%  '
